$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for row 2 (Receptor average / total expression value)
$ws.Range("M2").Value = 1.01111
$ws.Range("N2").Value = 3.03333

# Recomputed "Receptor derived specificity" columns (O,P,S,T) for rows 2-4
$ws.Range("O2").Value = 0.04063212692754557
$ws.Range("P2").Value = 0.04063212692754556
$ws.Range("S2").Value = 0.04063212692754557
$ws.Range("T2").Value = 0.04063212692754556

$ws.Range("O3").Value = 0.4065982422683317
$ws.Range("P3").Value = 0.4065982422683317
$ws.Range("S3").Value = 0.4065982422683317
$ws.Range("T3").Value = 0.4065982422683317

$ws.Range("O4").Value = 0.5527696308041227
$ws.Range("P4").Value = 0.5527696308041226
$ws.Range("S4").Value = 0.5527696308041227
$ws.Range("T4").Value = 0.5527696308041226

# Recomputed edge expression weights (Q,R) for row 2
$ws.Range("Q2").Value = 0.8422799077500001
$ws.Range("R2").Value = 7.580519169750001
